$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "RemoveInAndEx"

# ---- Row 1 (headers) ----
$ws.Range("A1").Value = "Execute"
$ws.Range("B1").Value = "Test Case"
$ws.Range("C1").Value = "Condition"
$ws.Range("D1").Value = "Expected"
$ws.Range("E1").Value = "ActualResult"
$ws.Range("F1").Value = "Result"
$ws.Range("G1").Value = "Revise"

# ---- Row 2 ----
$ws.Range("A2").Value = "Y"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "กดลบ"
$ws.Range("D2").Value = "ต้องการลบใช่หรือไม่ ?"

# ---- Row 3 ----
$ws.Range("A3").Value = "N"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "ไม่ใช่"
$ws.Range("D3").Value = "ข้อมูลขจะต้องอยู่ที่เดิม"

# ---- Row 4 ----
$ws.Range("A4").Value = "N"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "ใช่"
$ws.Range("D4").Value = "ข้อมูลของวันที่ทำรายการนั้นจะหายไป"

# ---- Column widths ----
# Target raw XML widths: D=26.09765625 (bestFit), E=16.296875 (bestFit, unchanged)
$ws.Columns("D:D").ColumnWidth = 25.333333333333332
$ws.Columns("E:E").ColumnWidth = 15.5

# ---- Alignment ----
# Style "center + vertical-center" on B1, E1:G1 (header cells that used to be centered/centered or left/centered)
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108
$ws.Range("E1:G1").HorizontalAlignment = -4108
$ws.Range("E1:G1").VerticalAlignment = -4108

# Style "center only" (horizontal center, default vertical) on A1, C1:D1 and columns A:B
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("C1:D1").HorizontalAlignment = -4108
$ws.Columns("A:B").HorizontalAlignment = -4108

# ---- Selection ----
$ws.Range("J1").Select() | Out-Null

# ---- Page setup ----
$ws.PageSetup.Orientation = 1
